$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) data rows keep their original Text storage,
# since many new price strings (e.g. "0.999", "1.00") would otherwise
# be auto-coerced to numbers by Excel, losing formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{Addr='D2'; Val='67.871.53'},
    @{Addr='E2'; Val='  -0.25%  '},
    @{Addr='D3'; Val='3.797.32'},
    @{Addr='E3'; Val='  -2.17%  '},
    @{Addr='D4'; Val='0.999'},
    @{Addr='E4'; Val='  -0.23%  '},
    @{Addr='D5'; Val='598.59'},
    @{Addr='E5'; Val='  -0.05%  '},
    @{Addr='D6'; Val='168.63'},
    @{Addr='E6'; Val='  -1.45%  '},
    @{Addr='D7'; Val='3.796.28'},
    @{Addr='E7'; Val='  -2.07%  '},
    @{Addr='E8'; Val='  -0.18%  '},
    @{Addr='D9'; Val='0.530'},
    @{Addr='E9'; Val='  +0.03%  '},
    @{Addr='E10'; Val='  +0.93%  '},
    @{Addr='E11'; Val='  +1.52%  '},
    @{Addr='E12'; Val='  +0.73%  '},
    @{Addr='E13'; Val='  +5.94%  '},
    @{Addr='D14'; Val='36.87'},
    @{Addr='E14'; Val='  -0.30%  '},
    @{Addr='D15'; Val='4.433.91'},
    @{Addr='E15'; Val='  -2.18%  '},
    @{Addr='D16'; Val='3.792.27'},
    @{Addr='E16'; Val='  -2.56%  '},
    @{Addr='D17'; Val='18.96'},
    @{Addr='E17'; Val='  +4.89%  '},
    @{Addr='D18'; Val='67.869.91'},
    @{Addr='E18'; Val='  -0.45%  '},
    @{Addr='D19'; Val='7.33'},
    @{Addr='E19'; Val='  -0.23%  '},
    @{Addr='E20'; Val='  +0.67%  '},
    @{Addr='D21'; Val='10.60'},
    @{Addr='E21'; Val='  -2.34%  '},
    @{Addr='D22'; Val='467.75'},
    @{Addr='E22'; Val='  +0.08%  '},
    @{Addr='D23'; Val='0.731'},
    @{Addr='E23'; Val='  -0.85%  '},
    @{Addr='D24'; Val='0.0000151'},
    @{Addr='E24'; Val='  -5.65%  '},
    @{Addr='D25'; Val='83.49'},
    @{Addr='E25'; Val='  +0.22%  '},
    @{Addr='E26'; Val='  +2.72%  '},
    @{Addr='D27'; Val='12.20'},
    @{Addr='E27'; Val='  +1.16%  '},
    @{Addr='D28'; Val='10.32'},
    @{Addr='E28'; Val='  +3.54%  '},
    @{Addr='D29'; Val='1.00'},
    @{Addr='E29'; Val='  +0.02%  '},
    @{Addr='E30'; Val='  -0.86%  '},
    @{Addr='D31'; Val='3.951.25'},
    @{Addr='E31'; Val='  -2.03%  '},
    @{Addr='D32'; Val='7.68'},
    @{Addr='E32'; Val='  -0.16%  '},
    @{Addr='D33'; Val='2.26'},
    @{Addr='E33'; Val='  -1.87%  '},
    @{Addr='D34'; Val='30.54'},
    @{Addr='E34'; Val='  -2.27%  '},
    @{Addr='E35'; Val='  -2.33%  '},
    @{Addr='D36'; Val='3.760.37'},
    @{Addr='E36'; Val='  -2.38%  '},
    @{Addr='E37'; Val='  +1.06%  '},
    @{Addr='D38'; Val='3.74'},
    @{Addr='E38'; Val='  -0.60%  '},
    @{Addr='D39'; Val='5.94'},
    @{Addr='E39'; Val='  +0.71%  '},
    @{Addr='E40'; Val='  -1.66%  '},
    @{Addr='E41'; Val='  -1.44%  '},
    @{Addr='D42'; Val='0.997'},
    @{Addr='E42'; Val='  -0.34%  '},
    @{Addr='E43'; Val='  +1.74%  '},
    @{Addr='B45'; Val='Cosmos'},
    @{Addr='C45'; Val='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Addr='D45'; Val='8.77'},
    @{Addr='E45'; Val='  +1.78%  '},
    @{Addr='B46'; Val='Stacks'},
    @{Addr='C46'; Val='https://coinranking.com/coin/mMPrMcB7+stacks-stx'},
    @{Addr='D46'; Val='1.97'},
    @{Addr='E46'; Val='  -0.55%  '},
    @{Addr='D47'; Val='408.00'},
    @{Addr='E47'; Val='  -3.70%  '},
    @{Addr='D48'; Val='46.30'},
    @{Addr='E48'; Val='  -1.92%  '},
    @{Addr='E49'; Val='  -8.06%  '},
    @{Addr='D50'; Val='142.07'},
    @{Addr='E50'; Val='  -0.70%  '},
    @{Addr='D51'; Val='0.0358'},
    @{Addr='E51'; Val='  +0.35%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}
